$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Tempo (ms)" column (B) for rows 2-11
$ws.Range("B2").Value = 2252
$ws.Range("B3").Value = 2559
$ws.Range("B4").Value = 2675
$ws.Range("B5").Value = 2731
$ws.Range("B6").Value = 2782
$ws.Range("B7").Value = 3132
$ws.Range("B8").Value = 3183
$ws.Range("B9").Value = 3417
$ws.Range("B10").Value = 3606
$ws.Range("B11").Value = 4093

# Update "Memória (KB)" column (C) for rows 5, 9, 10, 11
$ws.Range("C5").Value = 0.265625
$ws.Range("C9").Value = 0.28125
$ws.Range("C10").Value = 0.359375
$ws.Range("C11").Value = 249.109375

# Update summary rows (Média / Mediana). These cells hold text-formatted
# numbers (e.g. "3043.00"), not real numbers, so assigning them directly
# via .Value would coerce to a numeric cell and drop the formatting.
# Enter them as a text-literal formula, then paste-special as values so
# the result is a plain text cell (no residual formula, no style churn).
$ws.Range("B13").Formula = "=""3043.00"""
$ws.Range("C13").Formula = "=""25.17"""
$ws.Range("B14").Formula = "=""2957.00"""

$textCells = $ws.Range("B13:C14")
$textCells.Copy()
$textCells.PasteSpecial(-4163)
